# Snake.pptx - "Improvement Ideas" Folie hinzugefuegt.
#
# 1) Append a new slide at the end (position 9) using the "Titel und Inhalt"
#    (Title and Content) layout - the same layout used by every other
#    content slide in this deck - and fill in the title + bullet points.
# 2) Fix a small spelling typo ("histogramm" -> "histogram") on slide 5.

$p = $ppt.ActivePresentation

# --- 1) New slide: "Improvement ideas" -------------------------------------

$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$titleShape = $newSlide.Shapes.Item(1)
$titleShape.Name = "Titel 1"
$titleShape.TextFrame.TextRange.Text = "Improvement ideas"
$titleShape.TextFrame.TextRange.LanguageID = "en-GB"

$bodyShape = $newSlide.Shapes.Item(2)
$bodyShape.Name = "Inhaltsplatzhalter 2"
$bodyText = $bodyShape.TextFrame.TextRange

$bulletPoints = @(
    "Reduce number of states",
    "Smaller playground for the snake",
    "Using different base for the dictionary (not screen hash value)",
    "Make a different prime goal. Not eating the food but instead survive. "
)

$bodyText.Text = $bulletPoints[0]
$bodyText.LanguageID = "en-GB"
for ($i = 1; $i -lt $bulletPoints.Count; $i++) {
    $newPara = $bodyShape.TextFrame.TextRange.InsertAfter("`r" + $bulletPoints[$i])
    $newPara.LanguageID = "en-GB"
}

# --- 2) Typo fix on slide 5: "histogramm " -> "histogram " -----------------

$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $full = $shp.TextFrame.TextRange.Text
        $idx = $full.IndexOf("histogramm ")
        if ($idx -ge 0) {
            $fixRange = $shp.TextFrame.TextRange.Characters($idx + 1, 11)
            $fixRange.Text = "histogram "
        }
    }
}
